$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("N2").Value = 3.85
$ws.Range("P2").Value = 1.96
$ws.Range("Q2").Value = 2
# Row 3
$ws.Range("S3").Value = 3.15
$ws.Range("U3").Value = 2.32
# Row 4
$ws.Range("G4").Value = 6.4
$ws.Range("K4").Value = 3.9
$ws.Range("AH4").Value = 28
# Row 5
$ws.Range("S5").Value = 6.4
# Row 11
$ws.Range("G11").Value = 1.72
$ws.Range("R11").Value = 1.3
$ws.Range("AL11").Value = 42
$ws.Range("AO11").Value = 170
# Row 12
$ws.Range("G12").Value = 2.12
$ws.Range("I12").Value = 4.2
$ws.Range("J12").Value = 3.6
$ws.Range("S12").Value = 2.56
$ws.Range("W12").Value = 1.89
$ws.Range("AJ12").Value = 26
# Row 13
$ws.Range("J13").Value = 3.5
$ws.Range("R13").Value = 1.49
$ws.Range("AB13").Value = 13
# Row 15
$ws.Range("R15").Value = 1.42
$ws.Range("T15").Value = 1.58
# Row 16
$ws.Range("F16").Value = 2.62
$ws.Range("I16").Value = 2.86
$ws.Range("J16").Value = 3.25
$ws.Range("K16").Value = 110
$ws.Range("M16").Value = 1.06
$ws.Range("N16").Value = 1.1
$ws.Range("Q16").Value = 1.58
$ws.Range("R16").Value = 1.27
$ws.Range("S16").Value = 2.96
$ws.Range("T16").Value = 1.5
$ws.Range("V16").Value = 1.52
$ws.Range("X16").Value = 1000
$ws.Range("Y16").Value = 1000
$ws.Range("Z16").Value = 1000
$ws.Range("AA16").Value = 1000
$ws.Range("AB16").Value = 1000
$ws.Range("AC16").Value = 1000
$ws.Range("AD16").Value = 1000
$ws.Range("AE16").Value = 1000
$ws.Range("AF16").Value = 19
$ws.Range("AG16").Value = 1000
$ws.Range("AH16").Value = 1000
$ws.Range("AI16").Value = 1000
$ws.Range("AJ16").Value = 42
$ws.Range("AK16").Value = 30
$ws.Range("AL16").Value = 40
$ws.Range("AM16").Value = 1000
$ws.Range("AN16").Value = 1000
$ws.Range("AO16").Value = 1000
# Row 17
$ws.Range("G17").Value = 4.6
$ws.Range("H17").Value = 1.97
$ws.Range("Q17").Value = 1.64
# Row 18
$ws.Range("Q18").Value = 1.2
$ws.Range("T18").Value = 1.55
# Row 20
$ws.Range("P20").Value = 3.2
$ws.Range("Q20").Value = 1.31
$ws.Range("T20").Value = 1.58
$ws.Range("U20").Value = 1.97
# Row 21
$ws.Range("I21").Value = 3.3
$ws.Range("P21").Value = 2.56
$ws.Range("R21").Value = 1.63
$ws.Range("S21").Value = 2.24
$ws.Range("T21").Value = 1.52
$ws.Range("U21").Value = 2.56
$ws.Range("X21").Value = 28
$ws.Range("AJ21").Value = 32
$ws.Range("AK21").Value = 25
$ws.Range("AM21").Value = 60
$ws.Range("AO21").Value = 22
# Row 23
$ws.Range("S23").Value = 2.04
# Row 24
$ws.Range("I24").Value = 2.2
$ws.Range("J24").Value = 3.2
$ws.Range("V24").Value = 1.83
# Row 25
$ws.Range("F25").Value = 1.52
$ws.Range("G25").Value = 1.54
$ws.Range("H25").Value = 6.6
$ws.Range("I25").Value = 7
$ws.Range("K25").Value = 5.2
$ws.Range("N25").Value = 6
$ws.Range("Q25").Value = 1.57
$ws.Range("V25").Value = 1.16
$ws.Range("W25").Value = 2.84
$ws.Range("AA25").Value = 190
$ws.Range("AC25").Value = 11.5
$ws.Range("AF25").Value = 11
# Row 26
$ws.Range("K26").Value = 3.3
# Row 27
$ws.Range("F27").Value = 3.05
$ws.Range("G27").Value = 4.7
$ws.Range("H27").Value = 2.02
$ws.Range("I27").Value = 2.76
$ws.Range("K27").Value = 5.6
$ws.Range("N27").Value = 2.42
$ws.Range("S27").Value = 3.05
$ws.Range("V27").Value = 1.57
$ws.Range("W27").Value = 1.27
# Row 28
$ws.Range("G28").Value = 5.2
$ws.Range("I28").Value = 2.56
$ws.Range("J28").Value = 2.84
$ws.Range("N28").Value = 2.1
$ws.Range("V28").Value = 1.64
$ws.Range("W28").Value = 1.23
# Row 29
$ws.Range("G29").Value = 2.9
$ws.Range("H29").Value = 2.82
$ws.Range("K29").Value = 3.4
$ws.Range("AC29").Value = 7.6
$ws.Range("AH29").Value = 21
$ws.Range("AK29").Value = 36
# Row 30
$ws.Range("T30").Value = 1.69
# Row 31
$ws.Range("I31").Value = 4.5
$ws.Range("O31").Value = 1.48
$ws.Range("P31").Value = 1.66
$ws.Range("U31").Value = 1.84
$ws.Range("AD31").Value = 18.5
# Row 32
$ws.Range("J32").Value = 3.95
$ws.Range("AA32").Value = 21
$ws.Range("AL32").Value = 80
